$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# wdFindContinue = 1, wdReplaceAll = 2
$find.Execute("June 18, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "June 19, 2022", 2)
$find.Execute("August 17, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "August 18, 2022", 2)
